$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.401.36'
$ws.Range('E2').Value = '  +2.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.388.69'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.28'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.74'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('E9').Value = '  +5.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.592'
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.59'
$ws.Range('E11').Value = '  +2.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000282'
$ws.Range('E12').Value = '  +2.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '678.55'
$ws.Range('E13').Value = '  -4.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.63'
$ws.Range('E14').Value = '  +1.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.932.77'
$ws.Range('E15').Value = '  +1.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '69.475.09'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.120'
$ws.Range('E17').Value = '  +1.66%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.392.48'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('E19').Value = '  +0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.29'
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.905'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.20'
$ws.Range('E23').Value = '  +0.51%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '103.90'
$ws.Range('E24').Value = '  +3.39%  '
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.74'
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.18'
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.73'
$ws.Range('E29').Value = '  +1.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.03'
$ws.Range('E30').Value = '  -0.80%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.15'
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '557.74'
$ws.Range('E32').Value = '  -2.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.106'
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.55'
$ws.Range('E34').Value = '  +2.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.60'
$ws.Range('E35').Value = '  +1.57%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.691.14'
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.139'
$ws.Range('E38').Value = '  +4.36%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.43'
$ws.Range('E39').Value = '  +2.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.27'
$ws.Range('E40').Value = '  +2.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.70'
$ws.Range('E41').Value = '  +1.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0₃0699'
$ws.Range('E42').Value = '  +2.50%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.33'
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0424'
$ws.Range('E45').Value = '  +3.69%  '
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.42'
$ws.Range('E48').Value = '  +5.67%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.22'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('E51').Value = '  +2.37%  '
